# Updated cryptos list on Tue Sep 17 23:19:27 UTC 2024 with GitHub Actions
# Refreshes Sheet1 price/volume figures (columns D = Price, E = Volume(1h)).
# Numeric-looking Price values are written via NumberFormat "@" so Excel
# keeps them as text (matching the source data's inline-string storage),
# then the style is reset to Normal so no extra formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "60.094.71"
$ws.Range("E2").Value = "  +3.44%  "
$ws.Range("D3").Value = "2.331.75"
$ws.Range("E3").Value = "  +1.85%  "
$ws.Range("E4").Value = "  +0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "543.98"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.32%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "131.27"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("D9").Value = "2.328.59"
$ws.Range("E9").Value = "  +1.75%  "
$ws.Range("E10").Value = "  +1.29%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.50"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("E13").Value = "  +1.57%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "23.72"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("D15").Value = "2.745.04"
$ws.Range("E15").Value = "  +1.82%  "
$ws.Range("D16").Value = "60.053.17"
$ws.Range("E16").Value = "  +3.55%  "
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").Value = "2.331.71"
$ws.Range("E18").Value = "  +2.14%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "10.58"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.66%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "4.15"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.60%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.76"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +5.70%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "313.36"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.57%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.996"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.47%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "63.41"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.61%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.171"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.19%  "
$ws.Range("E26").Value = "  +0.21%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "7.87"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.71%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.35"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +7.11%  "
$ws.Range("E29").Value = "  +2.68%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "171.61"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("E31").Value = "  +13.35%  "
$ws.Range("D32").Value = "0.0₃0728"
$ws.Range("E32").Value = "  +1.52%  "
$ws.Range("E33").Value = "  +3.75%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.38"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +12.37%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.380"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.94%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "17.98"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.40%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  +0.00%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "4.13"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +6.14%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "319.15"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +11.10%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "38.02"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -1.08%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.52"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.07%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "140.25"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.11%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "3.45"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.23%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0944"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.37%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "19.41"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +7.71%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0496"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.54%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.559"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.20%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0212"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("D50").Value = "0.0₆0211"
$ws.Range("E50").Value = "  +15.25%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "11.02"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.77%  "
